$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Download All candidate names" - clear out the extra "Jejumar"/"Manuel" entries
# that had been appended to row 4 (A4 and D4), restoring it to only the
# name/nickname pair columns (B4, E4) like the other rows.
$ws.Range("A4").ClearContents()
$ws.Range("D4").ClearContents()

# UI update: move the active selection to E2.
$sel = $ws.Range("E2").Select()
